# "lock in current version"
# Replace the word / image / category table (rows 2-49) with the new
# set of words and image assignments, while keeping every row's
# category ("face" / "flower" / "none") consistent with its image.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    ,@("wenden", "none", "none")
    ,@("handeln", "face/face026.jpg", "face")
    ,@("ehren", "face/face004.jpg", "face")
    ,@("rufen", "none", "none")
    ,@("machen", "flower/flower021.jpg", "flower")
    ,@("segeln", "flower/flower032.jpg", "flower")
    ,@("enden", "none", "none")
    ,@("spielen", "face/face008.jpg", "face")
    ,@("zeugen", "face/face001.jpg", "face")
    ,@("opfern", "none", "none")
    ,@("mühen", "flower/flower028.jpg", "flower")
    ,@("wehtun", "flower/flower029.jpg", "flower")
    ,@("schalten", "none", "none")
    ,@("streichen", "face/face014.jpg", "face")
    ,@("husten", "face/face029.jpg", "face")
    ,@("klagen", "none", "none")
    ,@("mauern", "flower/flower013.jpg", "flower")
    ,@("drücken", "face/face022.jpg", "face")
    ,@("sparen", "none", "none")
    ,@("tragen", "face/face031.jpg", "face")
    ,@("gelten", "face/face002.jpg", "face")
    ,@("hören", "none", "none")
    ,@("grüßen", "face/face018.jpg", "face")
    ,@("treiben", "face/face030.jpg", "face")
    ,@("drohen", "none", "none")
    ,@("leugnen", "flower/flower008.jpg", "flower")
    ,@("stopfen", "flower/flower012.jpg", "flower")
    ,@("orten", "none", "none")
    ,@("tanzen", "face/face020.jpg", "face")
    ,@("reizen", "flower/flower001.jpg", "flower")
    ,@("weigern", "none", "none")
    ,@("helfen", "flower/flower026.jpg", "flower")
    ,@("buchen", "flower/flower027.jpg", "flower")
    ,@("ächzen", "none", "none")
    ,@("retten", "flower/flower025.jpg", "flower")
    ,@("füttern", "flower/flower031.jpg", "flower")
    ,@("kosten", "none", "none")
    ,@("leisten", "face/face011.jpg", "face")
    ,@("hoffen", "face/face010.jpg", "face")
    ,@("dauern", "none", "none")
    ,@("wagen", "flower/flower018.jpg", "flower")
    ,@("zögern", "flower/flower002.jpg", "flower")
    ,@("stören", "none", "none")
    ,@("antun", "flower/flower033.jpg", "flower")
    ,@("zählen", "face/face019.jpg", "face")
    ,@("bremsen", "none", "none")
    ,@("rühren", "flower/flower017.jpg", "flower")
    ,@("dringen", "face/face007.jpg", "face")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

